$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Alumno, Tarea, Fecha
$ws.Range("B9").Value = "Gaizka"
$ws.Range("C9").Value = "Footer con últimos empleados (extra)"

# Match the date format used by the cells above (D2:D8) before assigning the
# value, so Excel does not auto-create a brand new number format/style.
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("D9").Value = 45773

# Update the active selection to reflect the final cursor position after entry
$ws.Range("B10").Select()
